$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the data rows (2-17) to reflect the updated NATMI output
# (adds "ECs" as a 4th sending/target cluster, per Dr Hou's advice).

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Hgf"
$ws.Range("C2").Value = "Met"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 8.265822
$ws.Range("H2").Value = 24.797466
$ws.Range("I2").Value = 0.2082338764513023
$ws.Range("J2").Value = 0.2082338764513023
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.655851666666667
$ws.Range("N2").Value = 4.967555
$ws.Range("O2").Value = 0.03628213169899143
$ws.Range("P2").Value = 0.03628213169899143
$ws.Range("Q2").Value = 13.68697513507
$ws.Range("R2").Value = 123.18277621563
$ws.Range("S2").Value = 0.007555168929597659
$ws.Range("T2").Value = 0.007555168929597659

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Hgf"
$ws.Range("C3").Value = "Met"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 8.265822
$ws.Range("H3").Value = 24.797466
$ws.Range("I3").Value = 0.2082338764513023
$ws.Range("J3").Value = 0.2082338764513023
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.416382666666667
$ws.Range("N3").Value = 4.249148
$ws.Range("O3").Value = 0.03103501568568562
$ws.Range("P3").Value = 0.03103501568568562
$ws.Range("Q3").Value = 11.707567006552
$ws.Range("R3").Value = 105.368103058968
$ws.Range("S3").Value = 0.006462541621957286
$ws.Range("T3").Value = 0.006462541621957287

# Row 4: ECs -> M2
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Hgf"
$ws.Range("C4").Value = "Met"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 8.265822
$ws.Range("H4").Value = 24.797466
$ws.Range("I4").Value = 0.2082338764513023
$ws.Range("J4").Value = 0.2082338764513023
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.385314999999999
$ws.Range("N4").Value = 16.155945
$ws.Range("O4").Value = 0.1180001276707882
$ws.Range("P4").Value = 0.1180001276707882
$ws.Range("Q4").Value = 44.51405520392999
$ws.Range("R4").Value = 400.62649683537
$ws.Range("S4").Value = 0.0245716240066368
$ws.Range("T4").Value = 0.0245716240066368

# Row 5: ECs -> sCs
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Hgf"
$ws.Range("C5").Value = "Met"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 8.265822
$ws.Range("H5").Value = 24.797466
$ws.Range("I5").Value = 0.2082338764513023
$ws.Range("J5").Value = 0.2082338764513023
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 37.180664
$ws.Range("N5").Value = 111.541992
$ws.Range("O5").Value = 0.8146827249445348
$ws.Range("P5").Value = 0.8146827249445348
$ws.Range("Q5").Value = 307.328750465808
$ws.Range("R5").Value = 2765.958754192272
$ws.Range("S5").Value = 0.1696445418931105
$ws.Range("T5").Value = 0.1696445418931105

# Row 6: FAPs -> ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Hgf"
$ws.Range("C6").Value = "Met"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 11.28595333333333
$ws.Range("H6").Value = 33.85786
$ws.Range("I6").Value = 0.2843174958338682
$ws.Range("J6").Value = 0.2843174958338682
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.655851666666667
$ws.Range("N6").Value = 4.967555
$ws.Range("O6").Value = 0.03628213169899143
$ws.Range("P6").Value = 0.03628213169899143
$ws.Range("Q6").Value = 18.68786463692222
$ws.Range("R6").Value = 168.1907817323
$ws.Range("S6").Value = 0.01031564482817185
$ws.Range("T6").Value = 0.01031564482817185

# Row 7: FAPs -> FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Hgf"
$ws.Range("C7").Value = "Met"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 11.28595333333333
$ws.Range("H7").Value = 33.85786
$ws.Range("I7").Value = 0.2843174958338682
$ws.Range("J7").Value = 0.2843174958338682
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.416382666666667
$ws.Range("N7").Value = 4.249148
$ws.Range("O7").Value = 0.03103501568568562
$ws.Range("P7").Value = 0.03103501568568562
$ws.Range("Q7").Value = 15.98522867814222
$ws.Range("R7").Value = 143.86705810328
$ws.Range("S7").Value = 0.008823797942918956
$ws.Range("T7").Value = 0.008823797942918958

# Row 8: FAPs -> M2
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Hgf"
$ws.Range("C8").Value = "Met"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 11.28595333333333
$ws.Range("H8").Value = 33.85786
$ws.Range("I8").Value = 0.2843174958338682
$ws.Range("J8").Value = 0.2843174958338682
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.385314999999999
$ws.Range("N8").Value = 16.155945
$ws.Range("O8").Value = 0.1180001276707882
$ws.Range("P8").Value = 0.1180001276707882
$ws.Range("Q8").Value = 60.7784137753
$ws.Range("R8").Value = 547.0057239777
$ws.Range("S8").Value = 0.03354950080743523
$ws.Range("T8").Value = 0.03354950080743523

# Row 9: FAPs -> sCs
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Hgf"
$ws.Range("C9").Value = "Met"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 11.28595333333333
$ws.Range("H9").Value = 33.85786
$ws.Range("I9").Value = 0.2843174958338682
$ws.Range("J9").Value = 0.2843174958338682
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 37.180664
$ws.Range("N9").Value = 111.541992
$ws.Range("O9").Value = 0.8146827249445348
$ws.Range("P9").Value = 0.8146827249445348
$ws.Range("Q9").Value = 419.6192388063467
$ws.Range("R9").Value = 3776.57314925712
$ws.Range("S9").Value = 0.2316285522553422
$ws.Range("T9").Value = 0.2316285522553422

# Row 10: M2 -> ECs
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Hgf"
$ws.Range("C10").Value = "Met"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 19.51551966666667
$ws.Range("H10").Value = 58.546559
$ws.Range("I10").Value = 0.4916380138783083
$ws.Range("J10").Value = 0.4916380138783083
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.655851666666667
$ws.Range("N10").Value = 4.967555
$ws.Range("O10").Value = 0.03628213169899143
$ws.Range("P10").Value = 0.03628213169899143
$ws.Range("Q10").Value = 32.31480576591611
$ws.Range("R10").Value = 290.833251893245
$ws.Range("S10").Value = 0.01783767516776336
$ws.Range("T10").Value = 0.01783767516776336

# Row 11: M2 -> FAPs
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Hgf"
$ws.Range("C11").Value = "Met"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 19.51551966666667
$ws.Range("H11").Value = 58.546559
$ws.Range("I11").Value = 0.4916380138783083
$ws.Range("J11").Value = 0.4916380138783083
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.416382666666667
$ws.Range("N11").Value = 4.249148
$ws.Range("O11").Value = 0.03103501568568562
$ws.Range("P11").Value = 0.03103501568568562
$ws.Range("Q11").Value = 27.64144378685911
$ws.Range("R11").Value = 248.772994081732
$ws.Range("S11").Value = 0.01525799347239262
$ws.Range("T11").Value = 0.01525799347239262

# Row 12: M2 -> M2
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Hgf"
$ws.Range("C12").Value = "Met"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 19.51551966666667
$ws.Range("H12").Value = 58.546559
$ws.Range("I12").Value = 0.4916380138783083
$ws.Range("J12").Value = 0.4916380138783083
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 5.385314999999999
$ws.Range("N12").Value = 16.155945
$ws.Range("O12").Value = 0.1180001276707882
$ws.Range("P12").Value = 0.1180001276707882
$ws.Range("Q12").Value = 105.097220793695
$ws.Range("R12").Value = 945.874987143255
$ws.Range("S12").Value = 0.0580133484054531
$ws.Range("T12").Value = 0.05801334840545311

# Row 13: M2 -> sCs
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Hgf"
$ws.Range("C13").Value = "Met"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 19.51551966666667
$ws.Range("H13").Value = 58.546559
$ws.Range("I13").Value = 0.4916380138783083
$ws.Range("J13").Value = 0.4916380138783083
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 37.180664
$ws.Range("N13").Value = 111.541992
$ws.Range("O13").Value = 0.8146827249445348
$ws.Range("P13").Value = 0.8146827249445348
$ws.Range("Q13").Value = 725.5999795117253
$ws.Range("R13").Value = 6530.399815605528
$ws.Range("S13").Value = 0.4005289968326992
$ws.Range("T13").Value = 0.4005289968326992

# Row 14: sCs -> ECs
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Hgf"
$ws.Range("C14").Value = "Met"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.6276006666666666
$ws.Range("H14").Value = 1.882802
$ws.Range("I14").Value = 0.01581061383652123
$ws.Range("J14").Value = 0.01581061383652123
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1.655851666666667
$ws.Range("N14").Value = 4.967555
$ws.Range("O14").Value = 0.03628213169899143
$ws.Range("P14").Value = 0.03628213169899143
$ws.Range("Q14").Value = 1.039213609901111
$ws.Range("R14").Value = 9.35292248911
$ws.Range("S14").Value = 0.0005736427734585593
$ws.Range("T14").Value = 0.0005736427734585593

# Row 15: sCs -> FAPs
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Hgf"
$ws.Range("C15").Value = "Met"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.6276006666666666
$ws.Range("H15").Value = 1.882802
$ws.Range("I15").Value = 0.01581061383652123
$ws.Range("J15").Value = 0.01581061383652123
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 1.416382666666667
$ws.Range("N15").Value = 4.249148
$ws.Range("O15").Value = 0.03103501568568562
$ws.Range("P15").Value = 0.03103501568568562
$ws.Range("Q15").Value = 0.8889227058551109
$ws.Range("R15").Value = 8.000304352695998
$ws.Range("S15").Value = 0.0004906826484167544
$ws.Range("T15").Value = 0.0004906826484167544

# Row 16: sCs -> M2
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Hgf"
$ws.Range("C16").Value = "Met"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.6276006666666666
$ws.Range("H16").Value = 1.882802
$ws.Range("I16").Value = 0.01581061383652123
$ws.Range("J16").Value = 0.01581061383652123
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 5.385314999999999
$ws.Range("N16").Value = 16.155945
$ws.Range("O16").Value = 0.1180001276707882
$ws.Range("P16").Value = 0.1180001276707882
$ws.Range("Q16").Value = 3.379827284209999
$ws.Range("R16").Value = 30.41844555789
$ws.Range("S16").Value = 0.001865654451263035
$ws.Range("T16").Value = 0.001865654451263035

# Row 17: sCs -> sCs
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Hgf"
$ws.Range("C17").Value = "Met"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.6276006666666666
$ws.Range("H17").Value = 1.882802
$ws.Range("I17").Value = 0.01581061383652123
$ws.Range("J17").Value = 0.01581061383652123
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 37.180664
$ws.Range("N17").Value = 111.541992
$ws.Range("O17").Value = 0.8146827249445348
$ws.Range("P17").Value = 0.8146827249445348
$ws.Range("Q17").Value = 23.33460951350933
$ws.Range("R17").Value = 210.011485621584
$ws.Range("S17").Value = 0.01288063396338288
$ws.Range("T17").Value = 0.01288063396338288
